$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Object Code")

# Row 7 gets the new entry (ID 11). Set its label first so the string
# "PowerUp grabbed" is registered before "PowerUp draw" - this keeps the
# shared-string table ordering (and therefore the saved index numbers)
# aligned with the target workbook.
$ws.Range("A7").Value = 11
$ws.Range("B7").Value = "PowerUp grabbed"

# The old row 6 ("PowerUp") becomes "PowerUp draw".
$ws.Range("B6").Value = "PowerUp draw"

# Give column B an explicit custom width (matches the new column sizing
# added alongside the new row).
$ws.Columns.Item(2).ColumnWidth = 17.33
